# Refresh the crypto price/volume table (GitHub Actions data refresh).
# Cells in column D that look like plain numbers are written with a leading
# apostrophe so Excel keeps them as literal text (matching the source data,
# e.g. "13.00", "0.999") instead of silently parsing them as numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.833.57"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").Value = "3.635.10"
$ws.Range("E3").Value = "  +3.94%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'606.51"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").Value = "'199.02"
$ws.Range("E6").Value = "  +2.54%  "
$ws.Range("E7").Value = "  +0.83%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  +9.75%  "
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").Value = "'54.02"
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("E12").Value = "  +2.30%  "
$ws.Range("E13").Value = "  +1.18%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "4.209.29"
$ws.Range("E14").Value = "  +3.71%  "
$ws.Range("B15").Value = "BitcoinCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D15").Value = "'687.14"
$ws.Range("E15").Value = "  +16.10%  "
$ws.Range("D16").Value = "'13.00"
$ws.Range("E16").Value = "  +2.56%  "
$ws.Range("D17").Value = "70.885.93"
$ws.Range("E17").Value = "  +1.65%  "
$ws.Range("D18").Value = "3.629.49"
$ws.Range("E18").Value = "  +3.49%  "
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("E21").Value = "  +1.85%  "
$ws.Range("D22").Value = "'18.84"
$ws.Range("E22").Value = "  +3.32%  "
$ws.Range("D23").Value = "'5.41"
$ws.Range("E23").Value = "  +2.25%  "
$ws.Range("D24").Value = "'105.31"
$ws.Range("E24").Value = "  +4.04%  "
$ws.Range("D25").Value = "'4.63"
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("D26").Value = "'3.00"
$ws.Range("E27").Value = "  -2.85%  "
$ws.Range("D28").Value = "'9.86"
$ws.Range("E28").Value = "  +4.12%  "
$ws.Range("D29").Value = "'34.32"
$ws.Range("E29").Value = "  +3.80%  "
$ws.Range("D30").Value = "'4.61"
$ws.Range("E30").Value = "  +7.58%  "
$ws.Range("E31").Value = "  +2.55%  "
$ws.Range("D32").Value = "'12.19"
$ws.Range("E32").Value = "  -1.19%  "
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("D34").Value = "'63.34"
$ws.Range("E34").Value = "  +0.35%  "
$ws.Range("D35").Value = "0.0₃0871"
$ws.Range("E35").Value = "  +7.24%  "
$ws.Range("D36").Value = "3.951.10"
$ws.Range("E36").Value = "  +6.13%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E38").Value = "  -1.70%  "
$ws.Range("D39").Value = "'36.80"
$ws.Range("E39").Value = "  +1.78%  "
$ws.Range("D40").Value = "'505.14"
$ws.Range("E40").Value = "  +4.03%  "
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").Value = "'3.56"
$ws.Range("E42").Value = "  -2.71%  "
$ws.Range("D43").Value = "'0.137"
$ws.Range("E43").Value = "  +2.60%  "
$ws.Range("E44").Value = "  +10.46%  "
$ws.Range("D45").Value = "'0.0458"
$ws.Range("E45").Value = "  +1.66%  "
$ws.Range("D46").Value = "'3.52"
$ws.Range("E46").Value = "  +7.12%  "
$ws.Range("D47").Value = "'0.140"
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("E48").Value = "  +3.61%  "
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("D50").Value = "'0.000249"
$ws.Range("E50").Value = "  +1.91%  "
$ws.Range("E51").Value = "  +1.70%  "
